$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'35.332.26"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +0.75%  '
$ws.Range("D3").Value = "'1.898.90"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  +2.68%  '
$ws.Range("E4").Value = '  +0.30%  '
$ws.Range("D5").Value = "'243.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.75%  '
$ws.Range("D6").Value = "'0.651"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.92%  '
$ws.Range("E7").Value = '  +0.27%  '
$ws.Range("D8").Value = "'41.60"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.75%  '
$ws.Range("D9").Value = "'0.342"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +4.76%  '
$ws.Range("D10").Value = "'50.18"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +7.97%  '
$ws.Range("E11").Value = '  +2.82%  '
$ws.Range("E12").Value = '  +0.84%  '
$ws.Range("D13").Value = "'2.174.29"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +2.71%  '
$ws.Range("E14").Value = '  +6.35%  '
$ws.Range("D15").Value = "'0.693"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +2.81%  '
$ws.Range("D16").Value = "'1.898.64"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.71%  '
$ws.Range("E17").Value = '  +1.99%  '
$ws.Range("D18").Value = "'35.363.15"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("D19").Value = "'71.51"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.22%  '
$ws.Range("D20").Value = "'0.0₃0816"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +3.21%  '
$ws.Range("D21").Value = "'241.92"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +0.64%  '
$ws.Range("D22").Value = "'12.51"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.02%  '
$ws.Range("D23").Value = "'4.72"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.68%  '
$ws.Range("E24").Value = '  +0.28%  '
$ws.Range("E25").Value = '  +1.32%  '
$ws.Range("D26").Value = "'2.36"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +28.66%  '
$ws.Range("D27").Value = "'170.24"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +0.40%  '
$ws.Range("D28").Value = "'8.35"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +4.56%  '
$ws.Range("D29").Value = "'18.21"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +3.58%  '
$ws.Range("E30").Value = '  +2.35%  '
$ws.Range("E31").Value = '  +3.54%  '
$ws.Range("E33").Value = '  +0.35%  '
$ws.Range("D34").Value = "'0.933"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +19.23%  '
$ws.Range("D35").Value = "'4.11"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +2.47%  '
$ws.Range("D36").Value = "'1.73"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.86%  '
$ws.Range("D37").Value = "'2.03"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.69%  '
$ws.Range("E38").Value = '  +2.38%  '
$ws.Range("E39").Value = '  +4.05%  '
$ws.Range("E40").Value = '  +1.78%  '
$ws.Range("E41").Value = '  +14.50%  '
$ws.Range("D42").Value = "'15.87"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +6.44%  '
$ws.Range("D43").Value = "'89.43"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("D44").Value = "'1.338.84"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.29%  '
$ws.Range("E45").Value = '  +1.72%  '
$ws.Range("D46").Value = "'47.22"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +37.79%  '
$ws.Range("D47").Value = "'2.40"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -0.49%  '
$ws.Range("E48").Value = '  +1.81%  '
$ws.Range("D49").Value = "'12.25"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -9.80%  '
$ws.Range("D50").Value = "'6.50"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +0.47%  '
$ws.Range("D51").Value = "'2.084.23"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.22%  '
